$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: the instrument survey record was "rebased" - the inflation/ppp
# figures (H2:L2) no longer apply as-is; clear the stale computed values but
# keep cell H2's existing number-format style in place (blank out I2:L2
# entirely).
$ws.Range("H2:L2").ClearContents()

# Row 3: add the new wave ("w2") record for Cambodia covering 2020-2021.
$ws.Range("A3").Value = "CAS"
$ws.Range("B3").Value = "khm"
$ws.Range("C3").Value = "🇰🇭"
$ws.Range("D3").Value = "w2"
$ws.Range("E3").Value = "Cambodia"
$ws.Range("F3").Value = 2021
$ws.Range("G3").Value = "2020-2021"

# Match the saved selection/active cell state from the edit.
[void]$ws.Range("E4").Select()
